$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.710.09"
$ws.Range("E2").Value = "  +2.26%  "

$ws.Range("D3").Value = "2.422.49"
$ws.Range("E3").Value = "  +2.83%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.30%  "

$ws.Range("E9").Value = "  +4.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.52%  "

$ws.Range("E11").Value = "  +1.18%  "

$ws.Range("E12").Value = "  -2.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.24%  "

$ws.Range("D14").Value = "2.843.81"
$ws.Range("E14").Value = "  +2.48%  "

$ws.Range("D15").Value = "59.513.45"
$ws.Range("E15").Value = "  +2.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000140"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.33%  "

$ws.Range("D17").Value = "2.422.56"
$ws.Range("E17").Value = "  +3.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.69%  "

$ws.Range("E19").Value = "  +4.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("E21").Value = "  +3.31%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.65%  "

$ws.Range("E24").Value = "  +0.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.72%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "

$ws.Range("E28").Value = "  +6.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.62%  "

$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.41%  "

$ws.Range("E36").Value = "  +0.63%  "

$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.421"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "312.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.99%  "

$ws.Range("E42").Value = "  +2.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "143.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0967"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0524"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("E47").Value = "  +0.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.402"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.28%  "

$ws.Range("E49").Value = "  +2.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("E51").Value = "  +4.95%  "
